$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three new columns before the old "type" column (C), shifting the
#    old C/D/E/F (type/amount/date/userId) to F/G/H/I.
# ---------------------------------------------------------------------------
$ws.Range("C:E").Insert()

# ---------------------------------------------------------------------------
# 2. New header labels for the inserted columns.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "Unnamed: 0.1"
$ws.Cells.Item(1, 4).Value = "Unnamed: 0.1.1"
$ws.Cells.Item(1, 5).Value = "Unnamed: 0.1.1.1"

# ---------------------------------------------------------------------------
# 3. Full data grid (rows 2-20, columns A-I). $null means "leave blank".
#    Column order: row, A, B, C, D, E, F(type), G(amount), H(date), I(userId)
# ---------------------------------------------------------------------------
$rows = @(
  @(2, 0, 0, 0, 0, 0, "food", 10, "2021-08-29", "Dixon3220"),
  @(3, 1, 1, 1, 1, 1, "clothing", 1, "2021-07-26", "Dixon3221"),
  @(4, 2, 2, 2, 2, 2, "food", 20, "2021-06-14", "Dixon3220"),
  @(5, 3, 3, 3, 3, 3, "entertain", 80, "2021-08-20", "Dixon3220"),
  @(6, 4, 4, 4, 4, 4, "other", 422, "2021-06-14", "Dixon3220"),
  @(7, 5, 5, 5, 5, 5, "clothing", 60.9, "2021-08-21", "Dixon3220"),
  @(8, 6, 6, 6, 6, 6, "clothing", 20, "2021-08-30", "Dixon3220"),
  @(9, 7, 7, 7, 7, 7, "study", 20, "2021-08-30", "Dixon3220"),
  @(10, 8, 8, 8, 8, 8, "food", 4.8, "2021-08-30", "Dixon3220"),
  @(11, 9, 9, 9, 9, 9, "others", 20, "2021-08-30", "Dixon3220"),
  @(12, 10, 10, 10, 10, 10, "food", 40, "2021-08-30", "Dixon3220"),
  @(13, 11, 11, 11, 11, 11, "food", 40, "2021-08-30", "Dixon3220"),
  @(14, 12, 12, 12, 12, 12, "others", 20, "2021-08-30", "Dixon3220"),
  @(15, 13, 13, 13, 13, 13, "clothing", 30, "2021-09-01", "Dixon3220"),
  @(16, 14, 14, 14, 14, 14, "transport", 300, "2021-09-03", "Dixon3220"),
  @(17, 15, 15, 15, 15, 15, "clothing", 40, "2021-08-31", "Dixon3220"),
  @(18, 16, 16, 16, 16, $null, "entertain", 25, "2021-09-01", "Dixon3220"),
  @(19, 17, 17, $null, $null, $null, "study", 26, "2021-09-02", "Dixon3220"),
  @(20, 18, $null, $null, $null, $null, "clothing", 10, "2021-09-02", "Dixon3220")
)

# Columns that hold the text "date" values, which otherwise auto-coerce into
# serial date numbers when assigned through .Value - write them as explicit
# text (leading apostrophe) so they stay literal strings like "2021-08-29".
$dateCols = @(8)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($c = 1; $c -le 9; $c++) {
        $v = $r[$c]
        if ($null -eq $v) {
            continue
        }
        $cell = $ws.Cells.Item($rowNum, $c)
        if ($dateCols -contains $c) {
            $cell.Value = "'" + $v
        } else {
            $cell.Value = $v
        }
    }
}
